# "made all tests independent of each other"
#
# 1) Recolor the three custom fonts used on the sheet:
#      - team1/team2 labels (A2, C2) : red   -> green
#      - team3/team4 labels (A3, C3) : green -> black
#      - league3 label      (E4)     : blue  -> red
# 2) Give each "Point:" formula its own independent literal addend so the
#    cells no longer all evaluate to the same shared "0" result.
#
# Note: this runtime's PowerShell-COM surface has no RGB() helper, so the
# Font.Color values below are written directly as OLE_COLOR ints
# (0x00BBGGRR): 255 = red, 65280 = green, 0 = black.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Font color changes -------------------------------------------------
$ws.Range("A2").Font.Color = 65280   # was red,   now green (0,255,0)
$ws.Range("C2").Font.Color = 65280   # was red,   now green (0,255,0)

$ws.Range("A3").Font.Color = 0       # was green, now black (0,0,0)
$ws.Range("C3").Font.Color = 0       # was green, now black (0,0,0)

$ws.Range("E4").Font.Color = 255     # was blue,  now red   (255,0,0)

# --- Formula changes ------------------------------------------------------
$ws.Range("B2").Formula = "=0+20"
$ws.Range("D2").Formula = "=0+10"
$ws.Range("B3").Formula = "=0+30"
$ws.Range("D3").Formula = "=0+40"
